$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 14 date value (micro-adjustment of timestamp)
$ws.Cells.Item(14, 1).Value = 45814.39344563658

# Add new row 15 with latest price data
$ws.Cells.Item(15, 1).Value = 45815.39107168125
$ws.Cells.Item(15, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(15, 3).Value = "2Kg"
$ws.Cells.Item(15, 4).Value = "37,90€"

# Match the date formatting/style used by the other date cells in column A
$ws.Cells.Item(15, 1).NumberFormat = $ws.Cells.Item(14, 1).NumberFormat
